$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (column D) and Volume(1h) (column E) figures for the
# Sun Jan 22 09:38:38 UTC 2023 symbol-list refresh.
$updates = [ordered]@{
    "D2" = "302.03"
    "E2" = "-0.68%"
    "D3" = "37.41"
    "E3" = "7.25%"
    "D4" = "5.009"
    "E4" = "-3.97%"
    "D5" = "0.07824"
    "E5" = "0.20%"
    "D6" = "2.191"
    "E6" = "-7.84%"
    "D7" = "8.032"
    "E7" = "-0.20%"
    "E8" = "1.87%"
    "D9" = "0.9140"
    "E9" = "-2.16%"
    "D10" = "0.09735"
    "E10" = "-3.77%"
    "D11" = "0.1892"
    "E11" = "2.44%"
    "D12" = "0.08589"
    "E12" = "-0.34%"
    "D13" = "0.03524"
    "E13" = "6.34%"
    "D14" = "0.09965"
    "E14" = "0.70%"
    "D15" = "0.001482"
    "E15" = "-0.71%"
    "D16" = "0.005676"
    "E16" = "-1.37%"
    "D17" = "3.464"
    "E17" = "-0.26%"
    "D18" = "2.083"
    "E18" = "-3.02%"
    "E19" = "2.30%"
    "D21" = "4.765"
    "E21" = "10.66%"
    "E22" = "-0.97%"
    "D23" = "0.04642"
    "E23" = "1.49%"
    "D24" = "0.001230"
    "E24" = "1.06%"
    "D25" = "0.004798"
    "E25" = "8.10%"
    "E26" = "-7.56%"
    "E27" = "28.52%"
    "D39" = "0.01763"
    "E39" = "-0.72%"
    "D40" = "0.04726"
    "E40" = "-1.66%"
    "D41" = "0.008064"
    "E41" = "4.14%"
    "E42" = "-1.16%"
    "D43" = "0.007660"
    "E43" = "8.81%"
    "D44" = "0.002162"
    "E44" = "-5.92%"
    "D45" = "0.01041"
    "E45" = "9.94%"
    "D46" = "0.00006072"
    "E46" = "2.04%"
    "D47" = "0.00000000751"
    "E47" = "0.18%"
    "D48" = "7.640"
    "E48" = "179.60%"
    "D50" = "0.00002102"
    "E50" = "0.18%"
    "D51" = "0.0002002"
    "E51" = "0.18%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so numeric-looking strings (prices,
    # percentages) are stored as literal text instead of being
    # auto-converted to numbers by Excel's smart input parsing.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default (unstyled) look so the cell ends up exactly
    # like the other plain data cells in the sheet.
    $cell.Style = "Normal"
}
